## "wrapping up test file audit"
## The optimization_parameters sheet had a stray leftover row (A16="Sheet",
## B16=3, C16=4) that isn't part of the real parameter block - remove it.
## Removing it also drops the now-unused "Sheet" shared string and shifts
## every row below it up by one.

$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows.Item(16).Delete() | Out-Null

# Residual cursor position left on network_weights from the audit pass.
$weightsSheet = $wb.Worksheets.Item("network_weights")
$weightsSheet.Activate() | Out-Null
$weightsSheet.Range("D11").Select() | Out-Null

# Leave the workbook with "threshold_b" as the active/selected sheet, as in
# the saved file.
$thresholdSheet = $wb.Worksheets.Item("threshold_b")
$thresholdSheet.Activate() | Out-Null
$thresholdSheet.Select() | Out-Null
